$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings are created in this order by the original author:
#   32 = "NO", 33 = "beA/path1/ {Attrib}", 34 = "CRUD"
# Touch the cells in that same order so the regenerated sharedStrings.xml
# lines up with the target workbook.
$ws.Range("D2").Value = "NO"
$ws.Range("F6").Value = "beA/path1/ {Attrib}"
$ws.Range("C2").Value = "CRUD"

# Row 2 - Permission column CRUDL -> CRUD, KV2 Path Yes -> NO
$ws.Range("D3").Value = "NO"

# Row 3 - same change as row 2
$ws.Range("C3").Value = "CRUD"

# Row 4 - Permission column CRUDL -> CRUD only
$ws.Range("C4").Value = "CRUD"

# Row 5 - Permission column CRUDL -> CRUD only
$ws.Range("C5").Value = "CRUD"

# Row 6 - Result Failed -> Success, Permission CRUDL -> CRUD,
# Test Path beA/path1/Xyz -> beA/path1/ {Attrib}, and apply the
# green highlight fill (like row 4) across A6:F6 (including blank E6)
$ws.Range("A6").Value = "Success"
$ws.Range("C6").Value = "CRUD"
$ws.Range("A6:F6").Interior.Color = $ws.Range("A4").Interior.Color

# Scroll / selection change: top-left cell A7, active cell B7
$ws.Range("B7").Select()
$excel.ActiveWindow.ScrollRow = 7
